$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.310985565185547
$ws.Range("B1").Value = 3.186798810958862
$ws.Range("C1").Value = 2.552442312240601
$ws.Range("D1").Value = 2.4888014793396
$ws.Range("E1").Value = 2.20829439163208
